# Insert a new weekly "Ciboulette" record at row 316 of the data table,
# pushing the existing rows 316-344 down to 317-345 (dimension grows from
# A1:R344 to A1:R345).
#
# The new row reuses the same market/category/quality/unit/origin template
# as every other row in this block (A,B,C,E,F,G,H,I,K,L,M,N,O,P,Q,R are all
# identical across the block), only the date (D) and volume (J) differ per
# week: D316 = 44769 ("2022-07-27"), J316 = 120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 316:344 down to 317:345, leaving a blank (but formatted) row 316.
$ws.Rows(316).Insert()

# After the insert, row 317 holds what used to be row 316 - reuse it as the
# template for the new row's static columns (A:R), then overwrite D and J
# with the new record's own values.
$ws.Range("A316:R316").Value = $ws.Range("A317:R317").Value2
$ws.Range("D316").Value = 44769
$ws.Range("J316").Value = 120
